$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------------
# Row 1: new header labels (idx, idx2, Name, Date Start, Date End,
# (m3/s), (MW1), (MW2), (GWh) Winter, (GWh) Summer, (GWh) Year)
# ---------------------------------------------------------------------
$ws.Range("A1").Value = "idx"
$ws.Range("B1").Value = "idx2"
$ws.Range("C1").Value = "Name"
$ws.Range("D1").Value = "Date Start"
$ws.Range("E1").Value = "Date End"
$ws.Range("F1").Value = "(m3/s)"
$ws.Range("G1").Value = "(MW1)"
$ws.Range("H1").Value = "(MW2)"
$ws.Range("I1").Value = "(GWh) Winter"
$ws.Range("J1").Value = "(GWh) Summer"
$ws.Range("K1").Value = "(GWh) Year"

# ---------------------------------------------------------------------
# Table data: one row per plant, columns A-K.
# ---------------------------------------------------------------------
$data = @(
    ,@(1, '108600', 'Trümpler', '1816', '2000', '2.6', '0.42', '0.39', '0.6', '0.5', '1.1000000000000001')
    ,@(2, '108500', 'Freienstein', '1832', '2004', '12', '0.6', '0.6', '1.25', '1.2', '2.4500000000000002')
    ,@(3, '108300', 'Kollbrunn', '1832', '2010', '5.8', '0.37', '0.32', '0.8', '0.7', '1.5')
    ,@(4, '404900', 'Manegg', '1860', '1981', '4.9000000000000004', '0.64', '0.63', '0.9', '0.8', '1.7')
    ,@(5, '108400', 'Sennhof, Illnau-Effretikon', '1860', '2014', '5.7', '0.47', '0.47', '0.6', '0.6', '1.2')
    ,@(6, '405000', 'Letten', '1877', '2004', '100', '5', '4.26', '9.19', '12.44', '21.63')
    ,@(7, '404800', 'Waldhalde', '1895', '1967', '4.5', '1.55', '1.5', '4.21', '4.6500000000000004', '8.8699999999999992')
    ,@(8, '405100', 'Höngg', '1898', '1988', '50', '1.4', '1.3', '3.4', '3.6', '7')
    ,@(9, '108700', 'Eglisau', '1920', '2012', '500', '28.61', '32.450000000000003', '90.89', '103.09', '193.98')
    ,@(10, '404300', 'Pilgersteg', '1920', '2013', '1.5', '0.56000000000000005', '0.56000000000000005', '0.9', '1.1000000000000001', '2')
    ,@(11, '405200', 'Dietikon', '1933', $null, '100', '2.94', '2.7', '8.8000000000000007', '10.3', '19.100000000000001')
    ,@(12, '405300', 'Wettingen', '1933', '1964', '133', '4.8899999999999997', '4.63', '11.44', '14.09', '25.53')
    ,@(13, '404400', 'Etzelwerk Altendorf', '1937', '1992', '34', '54', '48.4', '56.4', '44.8', '101.2')
    ,@(14, '108800', 'Reckingen', '1941', '2004', '560', '6.79', '6.62', '20.98', '24.25', '45.24')
    ,@(15, '106400', 'Neuhausen', '1951', '2011', '29.9', '2.8', '2.4500000000000002', '10.35', '10.4', '20.75')
    ,@(16, '106500', 'Rheinau', '1956', '2005', '400', '19.72', '19.3', '42.29', '87.26', '129.55000000000001')
    ,@(17, '106200', 'Schaffhausen', '1964', $null, '500', '0.97', '0.85', '2.66', '3.16', '5.82')
    ,@(18, '108450', 'Pfungen', '1994', $null, '10', '0.35', '0.36', '0.55000000000000004', '0.45', '1')
    ,@(19, '405250', 'Wettingen-Dotierzentrale', '2007', $null, '12', '0.41', '0.37', '0.97', '1.68', '2.64')
    ,@(20, '108460', 'Hard Wülflingen', '2015', $null, '6.5', '0.57999999999999996', '0.57999999999999996', '1.25', '1.3', '2.5499999999999998')
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    if ($row[3] -ne $null) { $ws.Cells.Item($r, 4).Value = $row[3] }
    if ($row[4] -ne $null) { $ws.Cells.Item($r, 5).Value = $row[4] }
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
    $ws.Cells.Item($r, 9).Value = $row[8]
    $ws.Cells.Item($r, 10).Value = $row[9]
    $ws.Cells.Item($r, 11).Value = $row[10]
    $r = $r + 1
}

# ---------------------------------------------------------------------
# Clear the two rows that used to hold plant #19 and #20 (now blank,
# matching the filler rows further down the sheet).
# ---------------------------------------------------------------------
$ws.Range("A22:K23").ClearContents()

# ---------------------------------------------------------------------
# Selection, matching the authored edit.
# ---------------------------------------------------------------------
$ws.Range("A15:K15").Select()
